# DPLKKPS132_VerifKelPeserta.xlsx edit
# Update the "Kembalikan ke register" (reject) test data on sheet DPLKKPS132-002
# from register number M13220800000023 to M13220800000039, and move the
# active selection on sheet DPLKKPS132-001 to N2.

$wb = $excel.ActiveWorkbook

$wsReject = $wb.Worksheets.Item("DPLKKPS132-002")
$wsReject.Range("N2").Value = "M13220800000039"
$wsReject.Range("F2").Value = "Username : 31816;`nPassword : bni1234;`nRole : Penyelia Settlement;`nNo. Register : M13220800000039;`nStatus Verifikasi : 0 : Kembalikan ke Register;`nKeterangan Verifikasi : KEP.TRX.436 tidak disetujui"

$wsApprove = $wb.Worksheets.Item("DPLKKPS132-001")
$wsApprove.Activate()
$wsApprove.Range("N2").Select()
